$d = $word.ActiveDocument

$replacements = @(
    @{old="484÷8=60, 4"; new="358÷3=119, 1"},
    @{old="344÷7=49, 1"; new="238÷7=34, 0"},
    @{old="974÷7=139, 1"; new="646÷2=323, 0"},
    @{old="549÷4=137, 1"; new="715÷6=119, 1"},
    @{old="179÷9=19, 8"; new="298÷3=99, 1"},
    @{old="753÷5=150, 3"; new="337÷3=112, 1"},
    @{old="306÷6=51, 0"; new="807÷2=403, 1"},
    @{old="665÷8=83, 1"; new="730÷8=91, 2"},
    @{old="890÷5=178, 0"; new="756÷3=252, 0"},
    @{old="753÷8=94, 1"; new="961÷6=160, 1"},
    @{old="678÷7=96, 6"; new="565÷2=282, 1"},
    @{old="585÷6=97, 3"; new="346÷4=86, 2"},
    @{old="462÷4=115, 2"; new="820÷4=205, 0"},
    @{old="784÷4=196, 0"; new="104÷6=17, 2"},
    @{old="834÷3=278, 0"; new="964÷9=107, 1"},
    @{old="135÷8=16, 7"; new="566÷8=70, 6"},
    @{old="558÷8=69, 6"; new="997÷8=124, 5"},
    @{old="684÷2=342, 0"; new="390÷2=195, 0"},
    @{old="842÷8=105, 2"; new="932÷7=133, 1"},
    @{old="974÷4=243, 2"; new="937÷6=156, 1"},
    @{old="328÷8=41, 0"; new="656÷2=328, 0"},
    @{old="291÷9=32, 3"; new="173÷9=19, 2"},
    @{old="787÷2=393, 1"; new="299÷6=49, 5"},
    @{old="988÷4=247, 0"; new="209÷4=52, 1"},
    @{old="428÷8=53, 4"; new="429÷4=107, 1"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
